$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (FAPs -> FAPs, Artn/Gfra3)
$ws.Range("G2").Value = 0.6226999999999999
$ws.Range("I2").Value = 0.8812810914468889
$ws.Range("J2").Value = 0.881281091446889
$ws.Range("M2").Value = 0.939461
$ws.Range("N2").Value = 2.818383
$ws.Range("Q2").Value = 0.5850023646999999
$ws.Range("R2").Value = 5.265021282299999
$ws.Range("S2").Value = 0.8812810914468889
$ws.Range("T2").Value = 0.881281091446889

# Row 3 (MuSCs -> FAPs, Artn/Gfra3)
$ws.Range("E3").Value = 1
$ws.Range("F3").Value = 0.3333333333333333
$ws.Range("G3").Value = 0.04989433333333334
$ws.Range("H3").Value = 0.149683
$ws.Range("I3").Value = 0.07061334918422178
$ws.Range("J3").Value = 0.07061334918422178
$ws.Range("M3").Value = 0.939461
$ws.Range("N3").Value = 2.818383
$ws.Range("Q3").Value = 0.04687378028766667
$ws.Range("R3").Value = 0.421864022589
$ws.Range("S3").Value = 0.07061334918422178
$ws.Range("T3").Value = 0.07061334918422178

# Row 4 (Neutrophils -> FAPs, Artn/Gfra3)
$ws.Range("G4").Value = 0.03399066666666666
$ws.Range("I4").Value = 0.04810555936888933
$ws.Range("J4").Value = 0.04810555936888933
$ws.Range("M4").Value = 0.939461
$ws.Range("N4").Value = 2.818383
$ws.Range("Q4").Value = 0.03193290569733333
$ws.Range("R4").Value = 0.287396151276
$ws.Range("S4").Value = 0.04810555936888933
$ws.Range("T4").Value = 0.04810555936888933

$wb.Save()
